$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13 - this shifts rows 13..21 down to
# 14..22, carrying row heights/formatting with them exactly as needed.
$ws.Rows.Item(13).Insert()

# Row 13 (new) only has B/C populated (docente responsible), no A label.
# Copy number/alignment/font formatting from the (now-shifted) row 14 cells
# so the styles match exactly (wrap text / red font) instead of inheriting
# the bold "label" style from column A's column-level default.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$ws.Range("B13").Value = "5840514 - Graziela Zamponi"
$ws.Range("C13").Value = "5840514 - Graziela Zamponi"
$ws.Range("A13").Clear()

# Objetivos: full objectives text (row 10)
$objetivos = "Propiciar ao aluno o conhecimento dos gêneros por meio dos quais ele deverá agir linguisticamente no espaço acadêmico (Objetivo Geral); 2. Ler e redigir resumos acadêmicos e relatórios de pesquisa experimental, além de reconhecer as características de uma resenha (Objetivo Específico); 3. Dominar técnicas de escrita adequadas aos gêneros acadêmicos (Objetivo Específico)."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Programa resumido: (row 14 after shift)
$resumido = "O texto escrito da esfera acadêmica. Gêneros acadêmicos."
$ws.Range("B14").Value = $resumido
$ws.Range("C14").Value = $resumido

# Programa: (row 16 after shift)
$programa = "1. O texto escrito na esfera acadêmica`n    Aspectos constitutivos do texto escrito`n    Fatores de legibilidade`n    Coesão`n2. Gêneros acadêmicos`n    Noções de gêneros. Gêneros acadêmicos`n    Resumo e resenha`n    Relatório de pesquisa experimental"
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Método: (row 19 after shift) now holds the "N1 = Prova" criteria text
$metodo = "N 1  = Prova= 10,0`nN 2 = 1ª NP + 2ª NP  (ver abaixo)"
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Critério: (row 20 after shift) now holds the partial-grade breakdown text
$criterio = "1ª Nota Parcial - Resumo= 5,0`n2ª Nota Parcial - Relatório=5,0 NOTA FINAL = N1 + N2/ 2"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Norma de recuperação: (row 21 after shift)
$recuperacao = "Ao aluno que não alcançar a média 5,0 (cinco) no final do período letivo será dada uma recuperação, por meio de uma prova."
$ws.Range("B21").Value = $recuperacao
$ws.Range("C21").Value = $recuperacao

# Bibliografia: (row 22, newly created by the shift)
$bibliografia = "1. FARACO, Carlos Alberto; TEZZA, Cristóvão. Oficina de texto.  6 ed. Petrópolis: Vozes, 2008.`n2. ILARI, Rodolfo.  Introdução à Semântica: brincando com a gramática.  São Paulo: Contexto, 2001.`n3. ______.  Introdução ao estudo do léxico: brincando com as palavras.  São Paulo: Contexto, 2002.`n4. KLEIMAN, Ângela. Texto e leitor: aspectos cognitivos da leitura. 4.ed.  Campinas: Pontes, 1995.`n5. KOCH, Ingedore Villaça.  A coesão textual.  São Paulo: Contexto, 2001.`n6. LIBERATO, Yara; FULGÊNCIO, Lúcia.   É possível facilitar a leitura: um guia para escrever claro.  São Paulo: Contexto, 2007.`n7. MACHADO, A.R (coord.); LOUSADA, E.; ABREU-TARDELLI, L. S.  Resumo.  São Paulo: Parábola Editorial, 2004.`n8. ______.   Resenha.  São Paulo: Parábola Editorial, 2004.`n9. MARCUSCHI, Luiz Antônio.  Da fala para a escrita: atividades de retextualização.  São Paulo: Cortez, 2000.`n10. SERAFINI, Maria José.    Como escrever textos. 5.ed. São Paulo: Globo, 1992."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
